$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.676.50'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '2.030.26'
$ws.Range("E3").Value = '  -0.97%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.606'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.93'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.375'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0820'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.17%  '

$ws.Range("E11").Value = '  +0.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.61'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.01%  '

$ws.Range("D13").Value = '2.331.66'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.98'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.01%  '

$ws.Range("E15").Value = '  +1.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.20'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.39%  '

$ws.Range("D17").Value = '2.032.11'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").Value = '37.637.54'
$ws.Range("E18").Value = '  -1.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.49'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.86'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -7.20%  '

$ws.Range("E21").Value = '  -1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.25'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.23%  '

$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("E24").Value = '  -1.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.33%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.60'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.98%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.128'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.73'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.47%  '

$ws.Range("E30").Value = '  -2.45%  '

$ws.Range("E31").Value = '  +0.49%  '

$ws.Range("E32").Value = '  +8.76%  '

$ws.Range("E33").Value = '  -3.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0603'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.37%  '

$ws.Range("E35").Value = '  -1.82%  '

$ws.Range("E36").Value = '  +1.79%  '

$ws.Range("E37").Value = '  +1.37%  '

$ws.Range("E38").Value = '  +4.15%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.95'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +6.06%  '

$ws.Range("D41").Value = '1.537.85'
$ws.Range("E41").Value = '  +1.25%  '

$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '95.90'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.55%  '

$ws.Range("E44").Value = '  -2.18%  '

$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.10'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.05'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.37%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.64%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.97'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.00'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("D51").Value = '2.220.79'
$ws.Range("E51").Value = '  -0.87%  '
